# koreksi file excel
# - Rename "Memory (KB)" header (I4) to "Memory (MB)"
# - Replace static Avg Load Time values in F5/F11/F17 with AVERAGE formulas over
#   their respective H-column blocks (mirrors the existing F8/F14/F20 pattern)
# - Move the active selection to K11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I header: "Memory (KB)" -> "Memory (MB)"
$ws.Range("I4").Value2 = "Memory (MB)"

# Turn the "Avg Load Time (ms)" summary cells into live averages of their H column block
$ws.Range("F5").Formula = "=AVERAGE(H5:H7)"
$ws.Range("F11").Formula = "=AVERAGE(H11:H13)"
$ws.Range("F17").Formula = "=AVERAGE(H17:H19)"

# Move selection to K11
$ws.Range("K11").Select()
